$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells to remain plain text even though their new values
# look like numbers (Excel would otherwise auto-convert them to floats).
$textCells = @("D5", "D6", "D10", "D15", "D20", "D22", "D23", "D24", "D25", "D32", "D36", "D37", "D39", "D40", "D44", "D46", "D47", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.871.75'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '3.140.65'
$ws.Range("E3").Value = '  +1.28%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '530.28'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").Value = '138.68'
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.140.37'
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("E9").Value = '  +3.03%  '
$ws.Range("D10").Value = '7.22'
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("E12").Value = '  +3.60%  '
$ws.Range("D13").Value = '3.681.30'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").Value = '25.59'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").Value = '58.004.86'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").Value = '3.145.00'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").Value = '12.76'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").Value = '353.20'
$ws.Range("E22").Value = '  +4.96%  '
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '5.78'
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '68.76'
$ws.Range("E25").Value = '  +3.16%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("E30").Value = '  +4.98%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '6.18'
$ws.Range("E32").Value = '  -5.05%  '
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("E35").Value = '  -0.82%  '
$ws.Range("D36").Value = '5.01'
$ws.Range("E36").Value = '  +8.18%  '
$ws.Range("D37").Value = '157.86'
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("D39").Value = '26.46'
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("D40").Value = '1.26'
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  +6.95%  '
$ws.Range("E43").Value = '  +7.63%  '
$ws.Range("D44").Value = '0.706'
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("D45").Value = '3.184.09'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0273'
$ws.Range("E46").Value = '  +4.80%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = '36.63'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").Value = '2.343.05'
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = '6.06'
$ws.Range("E51").Value = '  +0.86%  '
